$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Sheet2" after the existing "Sheet1" and make it active.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "HUM103"
$ws2.Range("B1").Value = "HUM303"

# Data rows
$ws2.Range("A2").Value = 3
$ws2.Range("B2").Value = 1

$ws2.Range("A3").Value = 15
$ws2.Range("B3").Value = 2

$ws2.Range("A4").Value = 29
$ws2.Range("B4").Value = 4

$ws2.Range("A5").Value = 32
$ws2.Range("B5").Value = 5

$ws2.Range("A6").Value = 58
$ws2.Range("B6").Value = 43

$ws2.Range("B7").Value = 44

$ws2.Range("B8").Select()
